$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row (2..142) -> original row whose A/B values should land there.
$mapping = @(38,39,40,41,42,43,44,45,46,47,48,49,15,16,17,18,19,20,21,22,23,24,25,26,27,50,51,52,53,54,55,56,57,58,2,3,4,5,6,7,8,73,74,75,76,77,78,28,29,30,31,32,33,34,35,36,37,59,60,61,62,63,64,65,66,67,68,69,70,71,72,9,10,11,12,13,14,108,109,110,111,112,113,114,115,116,117,118,119,131,132,133,134,135,136,137,138,139,140,141,142,79,80,81,82,83,84,85,86,87,88,89,90,120,121,122,123,124,125,126,127,128,129,130,91,92,93,94,95,96,97,98,99,100,101,102,103,104,105,106,107)

# Snapshot the original English (col A) and Japanese (col B) values for rows 2..142
# before overwriting anything, since the mapping permutes rows in place.
$origA = @{}
$origB = @{}
for ($r = 2; $r -le 142; $r++) {
    $origA[$r] = $ws.Cells.Item($r, 1).Value2
    $origB[$r] = $ws.Cells.Item($r, 2).Value2
}

for ($i = 0; $i -lt $mapping.Length; $i++) {
    $newRow = $i + 2
    $srcRow = $mapping[$i]
    $ws.Cells.Item($newRow, 1).Value = $origA[$srcRow]
    $ws.Cells.Item($newRow, 2).Value = $origB[$srcRow]
}
